# BOM workbook update:
#  - Insert a new "KaiFeng Electronics" / "KF-08P" column pair before the
#    existing "Line total" column (which shifts from P to R).
#  - Correct the LED order: rows 3 and 5 had GREEN/YELLOW (and their
#    matching patterns) swapped - put them back the right way round.
#  - Fix the saved view (no frozen/scrolled topLeftCell, selection on F5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns at P, shifting the old "Line total" column (P) to R ---
$ws.Columns("P:Q").Insert()

# The insert copies O8's currency style into both P8 and Q8; only Q8 should
# carry that inherited formatting - P8 stays empty/unstyled.
$ws.Range("P8").Clear()

# Match the new columns' width to the rest of the L:O block
$ws.Range("P1:Q1").ColumnWidth = $ws.Range("O1").ColumnWidth

# --- New supplier header + part number ---
$ws.Range("P1:Q1").Merge()
$ws.Range("P1").Value = "KaiFeng Electronics"
$ws.Range("P9").Value = "KF-08P"
$ws.Range("Q15").Value = "Sub total"

# --- Correct LED documentation error: rows 3 & 5 had their Value/Pattern swapped ---
$ws.Range("B3").Value = "GREEN"
$ws.Range("F3").Value = "150060VS75000"
$ws.Range("B5").Value = "YELLOW"
$ws.Range("F5").Value = "150060YS75000"

# --- Restore the saved view ---
$ws.Range("F5").Select()
